# Insert three new vendor names into the alphabetically sorted "vendors"
# sheet, keeping the existing sort order intact:
#   - "LG"      before "Lost Vape"
#   - "Samsung" before "SMArt Mods"
#   - "Sony"    before "Squape"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vendors")
$ws.Activate() | Out-Null

# 1) Insert "LG" just above the row that currently holds "Lost Vape" (row 44).
$ws.Rows.Item(44).Insert()
$ws.Cells.Item(44, 1).Value = "LG"

# 2) Insert "Sony" just above the row that currently holds "Squape".
#    "Squape" was originally row 58; after the previous insert it is row 59.
$ws.Rows.Item(59).Insert()
$ws.Cells.Item(59, 1).Value = "Sony"

# 3) Insert "Samsung" just above the row that currently holds "SMArt Mods"
#    (still row 56 at this point, since the Sony insert happened below it).
$ws.Rows.Item(56).Insert()
$ws.Cells.Item(56, 1).Value = "Samsung"

# Leave the cursor/selection where the author ended up after the edit.
$ws.Range("A58").Select() | Out-Null
